$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: cardholder name / card number
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number that must stay TEXT (matches source which
# stored it as an inline string, not a number) -- force text format first so
# the long digit string isn't coerced into a floating point number.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 13.08.2025"

# Row 6: transaction 1
$ws.Range("B6").Value = "17.08."
$ws.Range("C6").Value = "18.08."
$ws.Range("D6").Value = "EBAY MKTPLC EU ILRCXR"
$ws.Range("E6").Value = "173,72-"

# Row 7: transaction 2
$ws.Range("B7").Value = "18.08."
$ws.Range("C7").Value = "19.08."
$ws.Range("D7").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 74921070"
$ws.Range("E7").Value = "83,07-"

# Row 8: transaction 3
$ws.Range("B8").Value = "22.08."
$ws.Range("C8").Value = "23.08."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "25,47-"

# Row 9: transaction 4 (previously blank row, now populated)
$ws.Range("B9").Value = "23.08."
$ws.Range("C9").Value = "24.08."
$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 47092717"
$ws.Range("E9").Value = "38,29-"
# E9 used to be an empty "center + wrap" cell (style 13); the now-populated
# amount cell should look like the other amount cells in the column
# (right aligned, no wrap, no forced vertical centering -- style 17).
$ws.Range("E9").HorizontalAlignment = -4152
$ws.Range("E9").VerticalAlignment = -4107
$ws.Range("E9").WrapText = $False

# Row 12: closing balance date/amount
$ws.Range("D12").Value = "KONTOSTAND AM 27.08.2025"
$ws.Range("E12").Value = "320,55-"

# Row 13: next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 03.09.2025"
